$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N (14th column), shifting TeamSize/Limit/StartLimit right
$ws.Columns.Item(14).Insert()

# The new column should inherit the width of the column that used to be here
# (now shifted to column O, originally "TeamSize", ~9.63 OOXML width units)
$ws.Columns.Item(14).ColumnWidth = 8.8

# Set header for new column N1
$ws.Range("N1").Value = "MinU18"

# Set value for new data cell N2
$ws.Range("N2").Value = 2

# Move the active selection to N2
$ws.Range("N2").Select() | Out-Null
